# AZ-104 Azure storage - Ahmed Wahid.pptx
# Commit: "Add files via upload"
#
# The underlying edit: on the "Objectives" slide (slideId 276), the
# SmartArt/SmartArt-diagram bullet that used to read
#   "Being able to choose the right storage and replication type to increase cost efficiency "
# was retyped to read
#   "Being able to choose the right storage to increase cost efficiency "
# (i.e. "and replication type to" was replaced by "to").
#
# It also looks like the author had that slide open/selected last, so the
# document-level "last slide viewed" tag was updated to point at slide 4
# ("Objectives", slideId 276).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Fix up the SmartArt text on the "Objectives" slide.
# ---------------------------------------------------------------------

$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 276) {
        $targetSlide = $p.Slides.Item($i)
        break
    }
}

if ($targetSlide -eq $null) {
    $targetSlide = $p.Slides.Item(4)
}

$oldText = "Being able to choose the right storage and replication type to increase cost efficiency "
$newText = "Being able to choose the right storage to increase cost efficiency "

for ($si = 1; $si -le $targetSlide.Shapes.Count; $si++) {
    $shp = $targetSlide.Shapes.Item($si)
    if ($shp.HasSmartArt) {
        $allNodes = $shp.SmartArt.AllNodes
        for ($ni = 1; $ni -le $allNodes.Count; $ni++) {
            $node = $allNodes.Item($ni)
            $tr = $node.TextFrame2.TextRange
            if ($tr.Text -eq $oldText) {
                $tr.Text = $newText
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Record the last slide viewed (slide 4, "Objectives", slideId 276).
# ---------------------------------------------------------------------

$p.Tags.Add("LASTSLIDEVIEWED", "276,4,Objectives")
